$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.844.67"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.630.65"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'215.65"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "'0.5122"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.2564"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.245"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.72"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").Value = "1.853.32"
$ws.Range("E14").Value = "  -1.00%  "
$ws.Range("D15").Value = "'0.5522"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "'63.65"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "0.0₅7589"
$ws.Range("D18").Value = "25.868.07"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "'194.66"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "'4.416"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'9.854"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'6.016"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'1.888"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "'142.17"
$ws.Range("E26").Value = "  +0.95%  "
$ws.Range("D27").Value = "'0.1257"
$ws.Range("E27").Value = "  +4.96%  "
$ws.Range("D28").Value = "'6.759"
$ws.Range("E28").Value = "  -0.83%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "'3.234"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'3.178"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'1.547"
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'2.372"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'0.8957"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'0.5525"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("D38").Value = "'2.537"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "1.114.70"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").Value = "'5.571"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").Value = "'0.7948"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").Value = "'97.62"
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").Value = "1.777.04"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -9.71%  "
$ws.Range("D47").Value = "'0.4433"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").Value = "'1.000"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'0.05135"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("D51").Value = "'7.560"
$ws.Range("E51").Value = "  +3.26%  "
